# This edit re-shuffles the observation rows (rows 3-23) of the "Artfynd"
# sheet: the record that used to live in one row now lives in a different
# row (as if the underlying record list had been re-sorted/re-fetched).
# For every target row, the full record (Id / Ost / Nord / "Publik
# kommentar", and whether the blank Alder-Stadium..Metod placeholder
# cells K:N exist) moves in from a specific source row.
#
# Strategy: snapshot every source row's current A/Q/R/AC values first (so
# later writes never clobber a value that is still needed as someone
# else's source), then write the new values into place. K:N are always
# empty placeholder cells; only rows 4 and 8 actually change whether
# those placeholder cells are present at all, so those two rows are
# patched up explicitly at the end (row 4 loses them, row 8 gains them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indexes used below.
$colA  = 1    # Id
$colQ  = 17   # Ost
$colR  = 18   # Nord
$colAC = 29   # Publik kommentar

# Rows 3-23 hold the observation records that get re-shuffled.
$dataRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23)

# target row -> source row (where its new data currently lives)
$mapping = @{
    3  = 7
    4  = 8
    5  = 3
    6  = 9
    7  = 4
    8  = 10
    9  = 11
    10 = 12
    11 = 5
    12 = 13
    13 = 14
    14 = 15
    15 = 16
    16 = 17
    17 = 18
    18 = 19
    19 = 20
    20 = 21
    21 = 22
    22 = 23
    23 = 6
}

# Snapshot every source row's current values before making any changes.
$snapA  = @{}
$snapQ  = @{}
$snapR  = @{}
$snapAC = @{}

foreach ($r in $dataRows) {
    $snapA[$r]  = $ws.Cells.Item($r, $colA).Value2
    $snapQ[$r]  = $ws.Cells.Item($r, $colQ).Value2
    $snapR[$r]  = $ws.Cells.Item($r, $colR).Value2
    $snapAC[$r] = $ws.Cells.Item($r, $colAC).Value2
}

# Now apply the new values, pulled from the recorded source row.
foreach ($r in $dataRows) {
    $src = $mapping[$r]

    $ws.Cells.Item($r, $colA).Value2 = $snapA[$src]
    $ws.Cells.Item($r, $colQ).Value2 = $snapQ[$src]
    $ws.Cells.Item($r, $colR).Value2 = $snapR[$src]

    if ($snapAC[$src] -eq $null) {
        $ws.Cells.Item($r, $colAC).ClearContents()
    } else {
        $ws.Cells.Item($r, $colAC).Value2 = $snapAC[$src]
    }
}

# Row 4 now represents what used to be row 8, which had no Alder-Stadium /
# Kon / Aktivitet / Metod placeholder cells (K:N) at all -> remove them.
$ws.Range("K4:N4").ClearContents()

# Row 8 now represents what used to be row 10, which did have blank K:N
# placeholder cells -> (re)create them as blank cells on row 8. Copying a
# blank K:N range from a row that still has them (row 9) reproduces an
# empty-but-present cell, which a plain blank-string assignment cannot do.
$ws.Range("K9:N9").Copy($ws.Range("K8:N8"))
